# Script 1 - atualização automática de dados
#
# Adds the newly published 01/07/2025 (Q3-2025) quarter to each of the
# three region blocks (Brasil, Nordeste, Sergipe) in the "g13.4" sheet.
# Each block keeps its data sorted chronologically, so the new quarter is
# inserted as a new row immediately after each block's previous last row
# (01/04/2025), pushing the following block(s) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$varName = "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência"

# Helper: write a full data row (Região, Variável, Trimestre, Valor) at a
# given row number without Excel re-interpreting the dd/mm/yyyy quarter
# label as a date serial number.
function Set-DataRow($row, $regiao, $trimestre, $valor) {
    $ws.Range("A$row").Value = $regiao
    $ws.Range("B$row").Value = $varName

    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = $trimestre
    $ws.Range("C$row").Style = "Normal"

    $ws.Range("D$row").Value = $valor
}

# --- Brasil block ends at row 27 (01/04/2025). Insert new row 28 for the
#     Brasil 01/07/2025 entry; this pushes the Nordeste + Sergipe blocks
#     down by one row each (Nordeste: 28-53 -> 29-54, Sergipe: 54-79 -> 55-80).
$ws.Rows.Item(28).Insert()
Set-DataRow 28 "Brasil" "01/07/2025" 94.43000000000001

# --- Nordeste block now ends at row 54 (01/04/2025, shifted from 53).
#     Insert new row 55 for the Nordeste 01/07/2025 entry; this pushes the
#     Sergipe block down by one more row (55-80 -> 56-81).
$ws.Rows.Item(55).Insert()
Set-DataRow 55 "Nordeste" "01/07/2025" 92.16

# --- Sergipe block now ends at row 81 (01/04/2025, shifted from 79).
#     Append the new Sergipe 01/07/2025 entry as row 82.
Set-DataRow 82 "Sergipe" "01/07/2025" 92.26000000000001

Write-Output "Added 01/07/2025 quarter for Brasil, Nordeste and Sergipe."
